$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 3
    "G2" = 13.36072166666666
    "H2" = 40.082165
    "I2" = 0.1827559288203559
    "J2" = 0.1827559288203559
    "K2" = 3
    "M2" = 2.761807333333334
    "N2" = 8.285422000000001
    "Q2" = 36.89973907762555
    "R2" = 332.09765169863
    "S2" = 0.1827559288203559
    "T2" = 0.1827559288203559

    "E3" = 3
    "G3" = 15.33382733333333
    "H3" = 46.001482
    "I3" = 0.2097452462965232
    "J3" = 0.2097452462965232
    "K3" = 3
    "M3" = 2.761807333333334
    "N3" = 8.285422000000001
    "Q3" = 42.34907677726711
    "R3" = 381.141690995404
    "S3" = 0.2097452462965232
    "T3" = 0.2097452462965232

    "E4" = 3
    "G4" = 41.04065866666667
    "H4" = 123.121976
    "I4" = 0.561378635162985
    "J4" = 0.561378635162985
    "K4" = 3
    "M4" = 2.761807333333334
    "N4" = 8.285422000000001
    "Q4" = 113.3463920704302
    "R4" = 1020.117528633872
    "S4" = 0.561378635162985
    "T4" = 0.561378635162985

    "E5" = 3
    "G5" = 3.371704666666667
    "H5" = 10.115114
    "I5" = 0.04612018972013576
    "J5" = 0.04612018972013576
    "K5" = 3
    "M5" = 2.761807333333334
    "N5" = 8.285422000000001
    "Q5" = 9.311998674234223
    "R5" = 83.807988068108
    "S5" = 0.04612018972013576
    "T5" = 0.04612018972013576
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
